$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0. Keep a pristine copy of the current "hyperlink cell" look (blue/underline
#    font + boxed border) on a scratch cell far outside the used range, so it
#    can be re-applied later no matter what Hyperlinks.Add() does to styles.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("Z99").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 1. Update existing row 2 (data row) values for columns A:I
# ---------------------------------------------------------------------------
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A2").Value = "yeimy@gmail.com"
$ws.Range("B2").Value = "Yei456"
$ws.Range("C2").Value = "Yei*76H"
$ws.Range("D2").Value = "Yeimy Lorena "
$ws.Range("E2").Value = "Santander"
$ws.Range("F2").Value = 3134777820
$ws.Range("G2").Value = "/users/"
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 1
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:yeimy@gmail.com")

# ---------------------------------------------------------------------------
# 2. Build the new "Update" mini-table headers (J1:O1) re-using the existing
#    header look (bold font / medium borders) but with a new fill colour.
# ---------------------------------------------------------------------------
# J1 reuses the "left edge" header style (same as A1 - all-sides border).
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# K1:O1 reuse the "inner" header style (same as B1 - no left border).
$ws.Range("B1").Copy()
$ws.Range("K1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-colour just the fill of the new header block (Gold, Accent 4, Lighter 80%).
$ws.Range("J1:O1").Interior.Color = 13431551

$ws.Range("J1").Value = "emailUpdate"
$ws.Range("K1").Value = "usernameUpdate"
$ws.Range("L1").Value = "passwordUpdate"
$ws.Range("M1").Value = "nameUpdate"
$ws.Range("N1").Value = "addressUpdate"
$ws.Range("O1").Value = "phoneUpdate"

# ---------------------------------------------------------------------------
# 3. New "Update" data row (J2:O2)
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = "Marina@gmail.com"
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:Marina@gmail.com")

$ws.Range("K2").Value = "Mari369"
$ws.Range("L2").Value = "Mar_3$"
$ws.Range("M2").Value = "Luz Marina Sosa"
$ws.Range("N2").Value = "Sucre/Sder"
$ws.Range("O2").Value = 321654987

# Give the new plain data cells (K2:O2) the same look as the rest of row 2.
$ws.Range("B2").Copy()
$ws.Range("K2:O2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K2").Value = "Mari369"
$ws.Range("L2").Value = "Mar_3$"
$ws.Range("M2").Value = "Luz Marina Sosa"
$ws.Range("N2").Value = "Sucre/Sder"
$ws.Range("O2").Value = 321654987

# ---------------------------------------------------------------------------
# 4. Re-apply the pristine hyperlink-cell look to A2 & J2 (Hyperlinks.Add
#    always stamps its own style, so restore the original box+font here) and
#    clean up the scratch cell.
# ---------------------------------------------------------------------------
$ws.Range("Z99").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z99").Clear()
$ws.Range("A2").Value = "yeimy@gmail.com"
$ws.Range("J2").Value = "Marina@gmail.com"

# ---------------------------------------------------------------------------
# 5. Cosmetics: widen the new columns and move the view over to show them.
# ---------------------------------------------------------------------------
$ws.Columns("J:O").AutoFit()
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("O7").Select()
